$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "예상 밖의 여정"
$ws.Range("E18").Value = "https://freesearch.pe.kr/archives/5272"

$ws.Range("D32").Value = "불균형 데이터 오버샘플링 기법: SMOTE, ADASYN, SMOTE-Tomek Link"
$ws.Range("E32").Value = "https://dodonam.tistory.com/509"
